$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "soldering status" values added in column N (rows 2-14).
# Rows 2-7 previously had no cell in column N at all; rows 8-14 had an
# empty (but styled) cell in column N. Writing .Value preserves the
# existing per-cell style and inherits the row's style for brand new cells.
$ws.Range("N2").Value  = "Solderinng"
$ws.Range("N3").Value  = "Desoldering"
$ws.Range("N4").Value  = "Soldering"
$ws.Range("N5").Value  = "No soldering"
$ws.Range("N6").Value  = "Soldering"
$ws.Range("N7").Value  = "Soldering"
$ws.Range("N8").Value  = "Desoldering"
$ws.Range("N9").Value  = "Soldering"
$ws.Range("N10").Value = "No Soldering"
$ws.Range("N11").Value = "No soldering"
$ws.Range("N12").Value = "Desoldering"
$ws.Range("N13").Value = "Soldering"
$ws.Range("N14").Value = "Soldering"

# Column H (8) got narrower -- was auto-sized to fit the long description
# text, now shrunk back down to a fixed width.
$ws.Columns.Item(8).ColumnWidth = 16.3

# Scroll/selection state: the view no longer has a frozen/offset
# topLeftCell, and the active selection moved to N19.
$ws.Range("N19").Select() | Out-Null
